$d = $word.ActiveDocument
$newText = "UNIVERSITY OF ABUJA`nDEPARTMENT OF PHYSICS`nFIRST SEMESTER 2018/2019 SESsSION EXAMINATION`nCOURSE TITLE: Basic Experimental Physies III`nCOURSE CODE: PHY 211`nINSTRUCTTONS: Answer question number 1 and any other one question.`n1. In an experiment perforned 5 times by a student to verify the length of his graph`nTIME ALLOWED: 1hr`nsheet, he recorded these observations.`nOBSERVATIONS LENGTH`n(CM)_`n31.33`n31.15:`n1.20`n51.02`n31.20`n(a). Determine the standard deviation`n(6). Determine the standard error`nCHence, what is the average paper length?`n2(a). The following values were obtained to determine acceleration due to gravity by`nusing bifilar suspension for 20 0scillations.`nLcm)_`n40.0`n50.0`nS/n`nTime t(S)`n|Period T (S)T)`n1D,0`nT. .0`n60.0`n70.0`n80.0`n19.0`n23.`ni.) Copy and complete the table above`nii.) Plot a graph of T on the vertical axis and L on the horizontal axis`nii). Determine the slope of the graph and state two precautions taken during the`nexperiment`nb. Define the term `"cóuple`" as it relates to rotational or oscillatory systems`n3a. What do you understand by precision and accuracy?`nb. Define each parameter in the straight line equetion: y= mx+c"
$d.Paragraphs.Item(1).Range.Text = $newText
